$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numbers formatted with "." as a thousands
# separator (e.g. "29.411.62", "1.001"), so they must stay literal TEXT --
# otherwise Excel auto-converts numeric-looking entries (like "1.001") into
# real numbers and mangles trailing zeros / separators. Temporarily mark the
# column as Text before writing, then restore the original (default) format
# and style so no visible formatting change is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.411.62"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.896.67"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "324.33"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.4763"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").Value = "0.4057"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").Value = "0.08023"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "1.002"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "23.32"
$ws.Range("E11").Value = "  +3.93%  "
$ws.Range("D12").Value = "1.864.90"
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("D13").Value = "5.928"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").Value = "7.059"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "89.56"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "0.06664"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "17.64"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "29.406.62"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "5.522"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "11.70"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "2.115.41"
$ws.Range("E25").Value = "  -3.27%  "
$ws.Range("D26").Value = "154.54"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "6.040"
$ws.Range("E28").Value = "  +4.65%  "
$ws.Range("D29").Value = "2.087"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").Value = "117.97"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "1.022"
$ws.Range("E31").Value = "  -5.50%  "
$ws.Range("D32").Value = "0.09458"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "1.388"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").Value = "3.530"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "5.363"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "0.06039"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").Value = "1.168"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").Value = "0.5856"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "7.802"
$ws.Range("E40").Value = "  -7.70%  "
$ws.Range("D41").Value = "0.1841"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "10.10"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "2.432"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.288"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("D45").Value = "0.07714"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").Value = "0.5496"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Value = "1.919"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "112.97"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "0.2964"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").Value = "43.61"
$ws.Range("E51").Value = "  -1.29%  "

# Restore original formatting/style on the Price column.
$priceRange.Style = "Normal"
